$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.759
$ws.Range("D5").Value = -7.907999999999999
$ws.Range("D6").Value = -7.831
$ws.Range("B11").Value = 6.305
$ws.Range("A12").Value = -21.489
$ws.Range("B23").Value = 8.73
$ws.Range("C24").Value = -12.572
$ws.Range("D27").Value = -7.795
$ws.Range("B28").Value = 5.142
$ws.Range("A32").Value = -21.319
$ws.Range("B32").Value = 7.694
$ws.Range("B34").Value = 7.292
$ws.Range("A36").Value = -20.724
$ws.Range("A38").Value = -20.311
$ws.Range("C38").Value = -11.85
$ws.Range("B42").Value = 9.074
$ws.Range("A46").Value = -21.728
$ws.Range("C52").Value = -11.826
$ws.Range("A54").Value = -21.37
$ws.Range("B54").Value = 5.406000000000001
$ws.Range("A55").Value = -22.311
$ws.Range("D55").Value = -7.951000000000001
$ws.Range("A67").Value = -21.473
$ws.Range("A69").Value = -21.422
$ws.Range("A72").Value = -21.722
$ws.Range("C78").Value = -12.484
$ws.Range("D80").Value = -7.935
$ws.Range("C83").Value = -13.593
$ws.Range("C85").Value = -12.335
$ws.Range("C86").Value = -13.725
$ws.Range("A91").Value = -20.92
$ws.Range("D95").Value = -7.628
$ws.Range("C96").Value = -11.412
$ws.Range("B97").Value = 5.219
$ws.Range("D98").Value = -7.758
$ws.Range("A99").Value = -21.282
$ws.Range("B99").Value = 5.670999999999999
$ws.Range("B101").Value = 5.225
$ws.Range("C103").Value = -12.872
$ws.Range("A104").Value = -21.437
